$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price / Volume(1h) columns) for rows 2-51.
# Values that look like plain numbers are apostrophe-prefixed so Excel
# keeps storing/display them as text (matching the original formatting,
# e.g. trailing zeros such as '1.000' or '81.00').
$updates = @(
    @{ Row = 2; D = '28.720.30'; E = '  +1.51%  ' }
    @{ Row = 3; D = '1.808.79'; E = '  -0.25%  ' }
    @{ Row = 4; D = '''1.002'; E = '  +0.28%  ' }
    @{ Row = 5; D = '''328.21'; E = '  -3.03%  ' }
    @{ Row = 6; D = '''1.000'; E = '  +0.47%  ' }
    @{ Row = 7; D = '''0.4391'; E = '  +0.51%  ' }
    @{ Row = 8; D = '''0.3768'; E = '  +6.64%  ' }
    @{ Row = 9; D = '''44.58'; E = '  -2.13%  ' }
    @{ Row = 10; D = '''0.07713'; E = '  +3.46%  ' }
    @{ Row = 11; D = '''1.141'; E = '  -1.43%  ' }
    @{ Row = 12; D = '''22.71'; E = '  -1.27%  ' }
    @{ Row = 13; D = '''1.000'; E = '  +0.28%  ' }
    @{ Row = 14; D = '''6.275'; E = '  -0.18%  ' }
    @{ Row = 15; D = '''7.547'; E = '  +3.34%  ' }
    @{ Row = 16; D = '1.806.49'; E = '  -0.52%  ' }
    @{ Row = 17; D = '''0.00001094'; E = '  +0.67%  ' }
    @{ Row = 18; D = '''0.06733'; E = '  +0.89%  ' }
    @{ Row = 19; D = '''81.00'; E = '  -1.35%  ' }
    @{ Row = 20; D = '''1.000'; E = '  +0.22%  ' }
    @{ Row = 21; D = '''17.68'; E = '  +2.01%  ' }
    @{ Row = 22; D = '''6.290'; E = '  -2.71%  ' }
    @{ Row = 23; D = '28.690.60'; E = '  +1.35%  ' }
    @{ Row = 24; D = '''11.78'; E = '  -2.67%  ' }
    @{ Row = 25; D = '''2.445'; E = '  +3.00%  ' }
    @{ Row = 26; D = '''20.58'; E = '  -0.80%  ' }
    @{ Row = 27; D = '''154.40'; E = '  -0.42%  ' }
    @{ Row = 28; D = '''2.368'; E = '  -4.20%  ' }
    @{ Row = 29; D = '2.013.95'; E = '  -0.51%  ' }
    @{ Row = 30; D = '''1.300'; E = '  -0.64%  ' }
    @{ Row = 31; D = '''131.29'; E = '  -1.05%  ' }
    @{ Row = 32; D = '''3.970'; E = '  -2.23%  ' }
    @{ Row = 33; D = '''5.831'; E = '  -2.41%  ' }
    @{ Row = 34; D = '''0.09211'; E = '  -1.52%  ' }
    @{ Row = 35; D = '''0.2238'; E = '  +3.30%  ' }
    @{ Row = 36; D = '''12.22'; E = '  -0.97%  ' }
    @{ Row = 37; D = '''0.06335'; E = '  +1.08%  ' }
    @{ Row = 38; D = '''5.222'; E = '  +0.09%  ' }
    @{ Row = 39; D = '''0.6624'; E = '  -2.54%  ' }
    @{ Row = 40; D = '''0.02316'; E = '  -2.54%  ' }
    @{ Row = 41; D = '''1.204'; E = '  -1.30%  ' }
    @{ Row = 42; D = '''8.085'; E = '  -2.02%  ' }
    @{ Row = 43; D = '''1.435'; E = '  -3.79%  ' }
    @{ Row = 44; D = '''0.9996'; E = '  +0.38%  ' }
    @{ Row = 45; D = '''13.93'; E = '  -0.34%  ' }
    @{ Row = 46; D = '''0.6085'; E = '  -1.45%  ' }
    @{ Row = 47; D = '''3.798'; E = '  -1.85%  ' }
    @{ Row = 48; D = '''128.05'; E = '  -1.76%  ' }
    @{ Row = 49; D = '''2.031'; E = '  -0.90%  ' }
    @{ Row = 50; D = '''0.07081'; E = '  -0.42%  ' }
    @{ Row = 51; D = '''1.146'; E = '  -2.51%  ' }
)

foreach ($u in $updates) {
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("E" + $u.Row).Value = $u.E
}
